# "Search dodan, dodan broj likeova u wall i newsfeed"
# (Search enabled, like-count added to wall & newsfeed responses)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# post/newsfeed (row 22) and post/wall (row 24) response JSON now include
# a "likesNumber" field alongside the existing "liked" boolean.
$likesText = '{ "data" : [postId : {"postId":id, "text" : text, "url":"url", "timestamp" : timestamp,"senderId":id, "senderName":name, "senderLastname":lastname, "senderPicture":url, "senderUsername":username, "senderEmail":email,  "recipientId":id, "recipientName":name, "recipientLastname":lastname, "recipientPicture":url, "recipientUsername":username, "recipientEmail":email, "liked": boolean, "likesNumber" : number}, secondPost : {secondPost}, .. nthPost : {nthPost}] , "error" : [] }'

$ws.Range("D22").Value = $likesText
$ws.Range("D24").Value = $likesText

# search (row 29): response JSON key renamed from "userId" to "id", and the
# endpoint's Status flips from Inactive to Active now that it's live.
$searchText = '{ "data" : [{"id":id, "name":name, "lastName":lastname, "profilePicture":url, "username":username}, {second user}, {third user}, …,{nth user}] , "error" : [] }'

$ws.Range("D29").Value = $searchText
$ws.Range("F29").Value = "Active"
